$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test-Cases")

# Rows 23 through 31: change "Approved/Rejected" (col I) from Rejected to Approved,
# and clear the "ReasonToReject" (col J) cell entirely.
for ($r = 23; $r -le 31; $r++) {
    $ws.Cells.Item($r, 9).Value = "Approved"
    $ws.Cells.Item($r, 10).ClearContents()
}

# Update the view state to match (scroll position + active selection).
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("H38").Select()
